$d = $word.ActiveDocument

# Locate the "Number of Biological Children" paragraph -- it's the one
# immediately before the block of three paragraphs (Child 1 / Child 2 /
# "Add more if applicable") that are being collapsed into a single
# {children_ages} placeholder.
$anchorPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Number of Biological Children*") {
        $anchorPara = $p
    }
}
Write-Host "Anchor paragraph: $($anchorPara.Range.Text)"

# Insert a brand-new paragraph right after the anchor. A freshly-minted
# paragraph/run carries no leftover run formatting, so setting its .Text
# below produces a plain <w:t> (no stray xml:space="preserve").
$anchorPara.Range.InsertParagraphAfter()
$newPara = $anchorPara.Next()
$newPara.Range.Text = "{children_ages}"
Write-Host "New paragraph: $($newPara.Range.Text)"

# Remove the three obsolete paragraphs that followed the anchor (Child 1,
# Child 2, and the "(Add more if applicable)" note). Re-fetch via .Next()
# each time so navigation stays correct as paragraphs are removed.
for ($i = 0; $i -lt 3; $i++) {
    $obsolete = $newPara.Next()
    Write-Host "Deleting: $($obsolete.Range.Text)"
    $obsolete.Range.Delete()
}
